$wb = $excel.ActiveWorkbook

# --- Sheet "library_content": add framework_min_score / framework_max_score rows ---
$wsLib = $wb.Worksheets.Item("library_content")

# fix NIST CSF score (library_version) from 1 to 2
$wsLib.Range("B2").Value = 2

# insert two new rows before the existing "tab" rows (old rows 14-15 -> 16-17)
$wsLib.Rows.Item(14).Insert()
$wsLib.Rows.Item(15).Insert()

$wsLib.Range("A14").Value = "framework_min_score"
$wsLib.Range("B14").Value = 1

$wsLib.Range("A15").Value = "framework_max_score"
$wsLib.Range("B15").Value = 4

# update the selection to match the new layout
$null = $wsLib.Range("B19").Select()

# --- Sheet "scores": selection only change (C3 -> B4) ---
$wsScores = $wb.Worksheets.Item("scores")
$null = $wsScores.Range("B4").Select()

# restore the originally active sheet ("library_content") as the active tab
$null = $wsLib.Activate()
$null = $wsLib.Range("B19").Select()
